$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Insert a new column before the old column D ("Vai indeksē?").
#    Old D -> E, old E -> F, old F -> G.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).Insert()

# ---------------------------------------------------------------------------
# 2. New column D header + data ("Vai mūsdienīgo?" Jā/Nē column)
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 4).Value = "Vai mūsdienīgo?"

$colD = @{
    2  = "Nē"
    3  = "Jā"
    4  = "Jā"
    5  = "Nē"
    6  = "Nē"
    7  = "Nē"
    8  = "Nē"
    9  = "Nē"
    10 = "Nē"
    11 = "Nē"
    12 = "Nē"
    13 = "Nē"
    15 = "Jā"
    16 = "Nē"
    17 = "Nē"
    19 = "Jā"
    21 = "Nē"
    22 = "Nē"
    23 = "Jā"
    25 = "Nē"
    27 = "Nē"
    29 = "Nē"
}
foreach ($r in $colD.Keys) {
    $ws.Cells.Item($r, 4).Value = $colD[$r]
}

# ---------------------------------------------------------------------------
# 3. Column E (formerly D, "Vai indeksē?") gets "Nē (Normunds)" notes added
#    for a number of rows that previously had no value there.
# ---------------------------------------------------------------------------
$rowsE = @(2,3,4,6,7,8,9,10,12,13,14,15,16,18,19,21,22,23,25,27)
foreach ($r in $rowsE) {
    $ws.Cells.Item($r, 5).Value = "Nē (Normunds)"
}

# ---------------------------------------------------------------------------
# 4. Row 27 also gets a note in column C ("Kaut kādi tukšumi?")
# ---------------------------------------------------------------------------
$ws.Cells.Item(27, 3).Value = "Kaut kādi tukšumi?"

# ---------------------------------------------------------------------------
# 5. Apply the wrap/vertical-top formatting used across the body of the
#    sheet to every new D/E cell so their style matches the rest of the
#    table (style index 2 in the original workbook: vertical=top, wrap=1).
#    Only touch cells that actually received a value above - formatting an
#    empty cell would materialize a spurious blank <c> in the XML.
# ---------------------------------------------------------------------------
foreach ($r in $colD.Keys) {
    $c = $ws.Cells.Item($r, 4)
    $c.VerticalAlignment = -4160
    $c.WrapText = $true
}

foreach ($r in $rowsE) {
    $c = $ws.Cells.Item($r, 5)
    $c.VerticalAlignment = -4160
    $c.WrapText = $true
}

$ws.Cells.Item(27, 3).VerticalAlignment = -4160
$ws.Cells.Item(27, 3).WrapText = $true

# ---------------------------------------------------------------------------
# 6. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 109.85546875
$ws.Columns.Item(4).ColumnWidth = 12.42578125
$ws.Columns.Item(5).ColumnWidth = 14.85546875
$ws.Columns.Item(6).ColumnWidth = 12.28515625

# ---------------------------------------------------------------------------
# 7. Row heights (custom)
# ---------------------------------------------------------------------------
$ws.Rows.Item(12).RowHeight = 16.5
$ws.Rows.Item(15).RowHeight = 152.25
$ws.Rows.Item(29).RowHeight = 318

# ---------------------------------------------------------------------------
# 8. Sheet view: freeze both row 1 and column A, with the active pane at
#    bottom-right positioned at D1, and matching per-pane selections.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("D1").Select()
